$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns (M:P) for detailed expense info, matching the style
# used by the rest of the header row (bold / bordered / centered).
$ws.Range("A1").Copy()
$ws.Range("M1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("M1").Value = "General & Administrative Expenses in USD millions"
$ws.Range("N1").Value = "Operating Expenses in USD millions"
$ws.Range("O1").Value = "Interest Expense in USD millions"
$ws.Range("P1").Value = "Depreciation, Amortization & Accretion in USD millions"

# All data-row cells on this sheet are stored as plain text, so force text
# formatting before assigning numeric-looking strings (keeps "$69.93",
# "28.3", "3,667", etc. as text instead of being coerced to numbers).
$dataCells = @("B2", "F2", "G2", "M2", "N2", "O2", "P2")
foreach ($addr in $dataCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated quote/market data for row 2
$ws.Range("B2").Value = "$69.93"
$ws.Range("C2").Value = "-3.25(4.44%) 1D"
$ws.Range("D2").Value = "$300.76B"
$ws.Range("F2").Value = "28.3"
$ws.Range("G2").Value = "12.1"

# New detail values for row 2
$ws.Range("M2").Value = "3,667"
$ws.Range("N2").Value = "4,026"
$ws.Range("O2").Value = "368"
$ws.Range("P2").Value = "290"
